$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.9817236666666668
$ws.Range("H2").Value = 2.945171
$ws.Range("I2").Value = 0.04688329326954743
$ws.Range("J2").Value = 0.04832841473263862
$ws.Range("N2").Value = 0.443301
$ws.Range("O2").Value = 0.2007197005285124
$ws.Range("P2").Value = 0.2255281176915293
$ws.Range("Q2").Value = 0.1450663610523334
$ws.Range("R2").Value = 1.305597249471
$ws.Range("S2").Value = 0.009410400584853979
$ws.Range("T2").Value = 0.01089941640566756
$ws.Range("G3").Value = 0.9817236666666668
$ws.Range("H3").Value = 2.945171
$ws.Range("I3").Value = 0.04688329326954743
$ws.Range("J3").Value = 0.04832841473263862
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2573206666666667
$ws.Range("N3").Value = 0.771962
$ws.Range("O3").Value = 0.3495322172956783
$ws.Range("P3").Value = 0.3927334627925232
$ws.Range("Q3").Value = 0.2526177883891111
$ws.Range("R3").Value = 2.273560095502
$ws.Range("S3").Value = 0.01638722145062846
$ws.Range("T3").Value = 0.01898018566922236
$ws.Range("G4").Value = 0.9817236666666668
$ws.Range("H4").Value = 2.945171
$ws.Range("I4").Value = 0.04688329326954743
$ws.Range("J4").Value = 0.04832841473263862
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2429445
$ws.Range("N4").Value = 0.485889
$ws.Range("O4").Value = 0.3300043127697603
$ws.Range("P4").Value = 0.2471946410610837
$ws.Range("Q4").Value = 0.2385043653365
$ws.Range("R4").Value = 1.431026192019
$ws.Range("S4").Value = 0.01547168897580013
$ws.Range("T4").Value = 0.01194652513288579
$ws.Range("G5").Value = 0.9817236666666668
$ws.Range("H5").Value = 2.945171
$ws.Range("I5").Value = 0.04688329326954743
$ws.Range("J5").Value = 0.04832841473263862
$ws.Range("M5").Value = 0.08815366666666667
$ws.Range("N5").Value = 0.264461
$ws.Range("O5").Value = 0.119743769406049
$ws.Range("P5").Value = 0.1345437784548637
$ws.Range("Q5").Value = 0.08654254087011112
$ws.Range("R5").Value = 0.778882867831
$ws.Range("S5").Value = 0.005613982258264854
$ws.Range("T5").Value = 0.006502287524862901
$ws.Range("G6").Value = 5.154927333333333
$ws.Range("H6").Value = 15.464782
$ws.Range("I6").Value = 0.2461792235003055
$ws.Range("J6").Value = 0.2537674037418691
$ws.Range("N6").Value = 0.443301
$ws.Range("O6").Value = 0.2007197005285124
$ws.Range("P6").Value = 0.2255281176915293
$ws.Range("Q6").Value = 0.7617281472646666
$ws.Range("R6").Value = 6.855553325382
$ws.Range("S6").Value = 0.04941302001732303
$ws.Range("T6").Value = 0.05723168489737009
$ws.Range("G7").Value = 5.154927333333333
$ws.Range("H7").Value = 15.464782
$ws.Range("I7").Value = 0.2461792235003055
$ws.Range("J7").Value = 0.2537674037418691
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2573206666666667
$ws.Range("N7").Value = 0.771962
$ws.Range("O7").Value = 0.3495322172956783
$ws.Range("P7").Value = 0.3927334627925232
$ws.Range("Q7").Value = 1.326469338031556
$ws.Range("R7").Value = 11.938224042284
$ws.Range("S7").Value = 0.08604756984219013
$ws.Range("T7").Value = 0.09966295121541258
$ws.Range("G8").Value = 5.154927333333333
$ws.Range("H8").Value = 15.464782
$ws.Range("I8").Value = 0.2461792235003055
$ws.Range("J8").Value = 0.2537674037418691
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2429445
$ws.Range("N8").Value = 0.485889
$ws.Range("O8").Value = 0.3300043127697603
$ws.Range("P8").Value = 0.2471946410610837
$ws.Range("Q8").Value = 1.252361243533
$ws.Range("R8").Value = 7.514167461198
$ws.Range("S8").Value = 0.08124020546941153
$ws.Range("T8").Value = 0.06272994228097446
$ws.Range("G9").Value = 5.154927333333333
$ws.Range("H9").Value = 15.464782
$ws.Range("I9").Value = 0.2461792235003055
$ws.Range("J9").Value = 0.2537674037418691
$ws.Range("M9").Value = 0.08815366666666667
$ws.Range("N9").Value = 0.264461
$ws.Range("O9").Value = 0.119743769406049
$ws.Range("P9").Value = 0.1345437784548637
$ws.Range("Q9").Value = 0.4544257458335555
$ws.Range("R9").Value = 4.089831712502
$ws.Range("S9").Value = 0.02947842817138076
$ws.Range("T9").Value = 0.03414282534811199
$ws.Range("G10").Value = 1.8784265
$ws.Range("H10").Value = 3.756853
$ws.Range("I10").Value = 0.08970632314876403
$ws.Range("J10").Value = 0.06164760887349412
$ws.Range("N10").Value = 0.443301
$ws.Range("O10").Value = 0.2007197005285124
$ws.Range("P10").Value = 0.2255281176915293
$ws.Range("Q10").Value = 0.2775694486255
$ws.Range("R10").Value = 1.665416691753
$ws.Range("S10").Value = 0.01800582631793387
$ws.Range("T10").Value = 0.01390326918942275
$ws.Range("G11").Value = 1.8784265
$ws.Range("H11").Value = 3.756853
$ws.Range("I11").Value = 0.08970632314876403
$ws.Range("J11").Value = 0.06164760887349412
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2573206666666667
$ws.Range("N11").Value = 0.771962
$ws.Range("O11").Value = 0.3495322172956783
$ws.Range("P11").Value = 0.3927334627925232
$ws.Range("Q11").Value = 0.4833579592643333
$ws.Range("R11").Value = 2.900147755586
$ws.Range("S11").Value = 0.03135525003563013
$ws.Range("T11").Value = 0.02421107890576643
$ws.Range("G12").Value = 1.8784265
$ws.Range("H12").Value = 3.756853
$ws.Range("I12").Value = 0.08970632314876403
$ws.Range("J12").Value = 0.06164760887349412
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2429445
$ws.Range("N12").Value = 0.485889
$ws.Range("O12").Value = 0.3300043127697603
$ws.Range("P12").Value = 0.2471946410610837
$ws.Range("Q12").Value = 0.45635338682925
$ws.Range("R12").Value = 1.825413547317
$ws.Range("S12").Value = 0.02960347352180992
$ws.Range("T12").Value = 0.01523895854775746
$ws.Range("G13").Value = 1.8784265
$ws.Range("H13").Value = 3.756853
$ws.Range("I13").Value = 0.08970632314876403
$ws.Range("J13").Value = 0.06164760887349412
$ws.Range("M13").Value = 0.08815366666666667
$ws.Range("N13").Value = 0.264461
$ws.Range("O13").Value = 0.119743769406049
$ws.Range("P13").Value = 0.1345437784548637
$ws.Range("Q13").Value = 0.1655901835388333
$ws.Range("R13").Value = 0.9935411012329999
$ws.Range("S13").Value = 0.01074177327339011
$ws.Range("T13").Value = 0.008294302230547482
$ws.Range("G14").Value = 12.924656
$ws.Range("H14").Value = 38.773968
$ws.Range("I14").Value = 0.617231160081383
$ws.Range("J14").Value = 0.6362565726519981
$ws.Range("N14").Value = 0.443301
$ws.Range("O14").Value = 0.2007197005285124
$ws.Range("P14").Value = 0.2255281176915293
$ws.Range("Q14").Value = 1.909837643152
$ws.Range("R14").Value = 17.188538788368
$ws.Range("S14").Value = 0.1238904536084015
$ws.Range("T14").Value = 0.1434937471990689
$ws.Range("G15").Value = 12.924656
$ws.Range("H15").Value = 38.773968
$ws.Range("I15").Value = 0.617231160081383
$ws.Range("J15").Value = 0.6362565726519981
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.2573206666666667
$ws.Range("N15").Value = 0.771962
$ws.Range("O15").Value = 0.3495322172956783
$ws.Range("P15").Value = 0.3927334627925232
$ws.Range("Q15").Value = 3.325781098357333
$ws.Range("R15").Value = 29.932029885216
$ws.Range("S15").Value = 0.2157421759672296
$ws.Range("T15").Value = 0.2498792470021219
$ws.Range("G16").Value = 12.924656
$ws.Range("H16").Value = 38.773968
$ws.Range("I16").Value = 0.617231160081383
$ws.Range("J16").Value = 0.6362565726519981
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2429445
$ws.Range("N16").Value = 0.485889
$ws.Range("O16").Value = 0.3300043127697603
$ws.Range("P16").Value = 0.2471946410610837
$ws.Range("Q16").Value = 3.139974089592
$ws.Range("R16").Value = 18.839844537552
$ws.Range("S16").Value = 0.2036889448027387
$ws.Range("T16").Value = 0.157279215099466
$ws.Range("G17").Value = 12.924656
$ws.Range("H17").Value = 38.773968
$ws.Range("I17").Value = 0.617231160081383
$ws.Range("J17").Value = 0.6362565726519981
$ws.Range("M17").Value = 0.08815366666666667
$ws.Range("N17").Value = 0.264461
$ws.Range("O17").Value = 0.119743769406049
$ws.Range("P17").Value = 0.1345437784548637
$ws.Range("Q17").Value = 1.139355816805333
$ws.Range("R17").Value = 10.254202351248
$ws.Range("S17").Value = 0.07390958570301322
$ws.Range("T17").Value = 0.08560436335134132
